# Chap04: Correction of picture and Biblio
#
# 1) Refresh the cached "datetimeFigureOut" footer text (31/10/2017 -> 06/11/2017)
#    on the slide master and every slide layout.
# 2) Re-size/re-position the two "1.43 eV" / "1.36 eV" labels and correct their
#    numeric values (1.43 -> 1.2, 1.36 -> 1.0).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Date placeholders (slide master + every custom layout)
# ---------------------------------------------------------------------------
function Update-DatePlaceholders($shapes, [string]$text) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDatePh = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                $isDatePh = $true
            }
        } catch {
        }
        if ($isDatePh) {
            $shp.TextFrame.TextRange.Text = $text
        }
    }
}

$design = $p.Designs.Item(1)
$master = $design.SlideMaster

Update-DatePlaceholders $master.Shapes "06/11/2017"

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholders $layout.Shapes "06/11/2017"
}

# ---------------------------------------------------------------------------
# 2) "1.43 eV" -> "1.2 eV" and "1.36 eV" -> "1.0 eV" text boxes
# ---------------------------------------------------------------------------
# Points <-> EMU: 1 pt = 12700 EMU. Position/size setters take points as a
# (single precision) float, so nudge by a hair over the exact value to avoid
# the EMU landing one unit short after the float32 round-trip/truncation.
function EmuToPt([double]$emu) {
    return ($emu / 12700.0) + 0.00001
}

$s = $p.Slides.Item(1)

# Shape "ZoneTexte 124" (id 125) - the "1.43 eV" label
$shLabel1 = $s.Shapes.Item(13)
$shLabel1.Left   = EmuToPt 3042000
$shLabel1.Top    = EmuToPt 2044800
$shLabel1.Width  = EmuToPt 800219
$shLabel1.Height = EmuToPt 369332

$run1 = $shLabel1.TextFrame.TextRange.Characters(1, 5)
$run1.Text = "1.2 "

# Shape "ZoneTexte 126" (id 127) - the "1.36 eV" label
$shLabel2 = $s.Shapes.Item(15)
$shLabel2.Left   = EmuToPt 3042000
$shLabel2.Top    = EmuToPt 884255
$shLabel2.Width  = EmuToPt 800219
$shLabel2.Height = EmuToPt 369332

$run2 = $shLabel2.TextFrame.TextRange.Characters(1, 5)
$run2.Text = "1.0 "
